$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("O1").Value = "teamid"

$teamIds = @(
    "7c925ac6-1512-497a-8807-19aee53531d2",
    "c67f321d-fab0-4e2a-bdaf-5801469d37b5",
    "5b8d0d04-7fc0-4b95-ba05-d124fa4334ac",
    "d88ef14f-d2dc-4721-aacc-1f632f7347eb",
    "1c104287-7c20-42e4-9fbb-5a9141b323ea",
    "7508a993-6eee-47d8-a537-5d5aec6bb902",
    "eaa39d6b-eb3b-4ee5-bab1-1908c7aabb73",
    "2af463cf-4130-4f2e-9df8-d91f58f98780",
    "05d44392-c3eb-4e0e-a3df-28a547d2b0d3",
    "7c45e8da-a16c-4a6e-afb9-992af41f1f2b",
    "2ba5bbb5-70f8-46c3-b473-907be766d907",
    "037f0493-37c7-49c2-a571-bb9a1c67c0e9",
    "fe67fd6b-7f29-448c-a4ef-241e6f3e6b43",
    "97764f5d-102d-44f7-b018-515714df51bf",
    "8130a434-1d36-4aae-b6ad-1aa871f06165",
    "b00db273-1041-43c3-b471-4e86d8982b93",
    "b9af5cac-c771-40c1-8e03-8c46bd84538d",
    "1da80034-143c-4b79-b58a-96a68a75c340",
    "b623377a-c762-48a0-9124-9bb8a345dae6",
    "8107f4a7-7bb6-4207-940c-f34d2e82a8f0",
    "c9986db9-5430-4be4-8bbd-d816e112ac1c",
    "82f70bed-ebcd-4751-8245-7221838bc31d",
    "f1657738-70c3-46c1-abc4-b8a6a6832705"
)

for ($i = 0; $i -lt $teamIds.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 15).Value = $teamIds[$i]
}
